# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5's table (2nd shape) switches to a different built-in table
#    style (tableStyleId GUID changes, bandRow/firstRow stay as-is).
# 2) The deck's theme (theme1.xml, used by the one slide master / all
#    slides) swaps its 12-colour scheme from the "Red Violet" / Integral
#    palette to the default "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{E360AEE3-61A1-4D55-805A-32526AA6FC69}", $true)

# --- 2. Theme colour scheme (Integral/Red Violet -> Office) --------------
$slide1 = $p.Slides.Item(1)
$colors = $slide1.ThemeColorScheme
$colors.Colors(1).RGB  = 0         # dk1      000000
$colors.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388   # dk2      44546A
$colors.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501   # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407     # accent4  FFC000
$colors.Colors(9).RGB  = 12874308  # accent5  4472C4
$colors.Colors(10).RGB = 4697456   # accent6  70AD47
$colors.Colors(11).RGB = 12673797  # hlink    0563C1
$colors.Colors(12).RGB = 7491477   # folHlink 954F72
